$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3234.4849
$ws.Range("I33").Value = 4356.0835
$ws.Range("K33").Value = 4356.0835
$ws.Range("M33").Value = -4127.0835

$ws.Range("H64").Value = 4301.1313
$ws.Range("I64").Value = 3348.0667
$ws.Range("J64").Value = 7875.125
$ws.Range("K64").Value = 3348.0667
$ws.Range("L64").Value = 7875.125
$ws.Range("M64").Value = -3100.0667
$ws.Range("N64").Value = -8371.125

$ws.Range("H67").Value = 4301.1313
$ws.Range("I67").Value = 3348.0667
$ws.Range("J67").Value = 7875.125
$ws.Range("K67").Value = 3348.0667
$ws.Range("L67").Value = 7875.125
$ws.Range("M67").Value = -2490.0667
$ws.Range("N67").Value = -9591.125

$ws.Range("H74").Value = 4632
$ws.Range("I74").Value = 4500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3564

$ws.Range("H77").Value = 4632
$ws.Range("I77").Value = 4500
$ws.Range("K77").Value = 22500
$ws.Range("M77").Value = -17820

$ws.Range("H86").Value = 101200.75
$ws.Range("I86").Value = 1601
$ws.Range("K86").Value = 1601
$ws.Range("M86").Value = -478

$ws.Range("H89").Value = 101200.75
$ws.Range("I89").Value = 1601
$ws.Range("K89").Value = 8005
$ws.Range("M89").Value = -2389

$ws.Range("H132").Value = 1700.4166
$ws.Range("I132").Value = 1490.8182
$ws.Range("K132").Value = 4472.4546
$ws.Range("M132").Value = -1942.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3740.441
$ws.Range("I32").Value = 2471.5732
$ws.Range("J32").Value = 13199.272
$ws.Range("K32").Value = 2471.5732
$ws.Range("L32").Value = 13199.272
$ws.Range("M32").Value = -2184.5732
$ws.Range("N32").Value = -13773.272

$ws.Range("H45").Value = 2070.25
$ws.Range("I45").Value = 2187.3333
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2187.3333
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1810.3333
$ws.Range("N45").Value = -2754

$ws.Range("H61").Value = 1925.1765
$ws.Range("I61").Value = 2565.375
$ws.Range("J61").Value = 1356.1111
$ws.Range("K61").Value = 2565.375
$ws.Range("L61").Value = 1356.1111
$ws.Range("M61").Value = -2353.375
$ws.Range("N61").Value = -1780.1111

$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 30000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 30000
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -29314

$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 30000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 150000
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -146568

$ws.Range("H74").Value = 2913
$ws.Range("I74").Value = 3043.4285
$ws.Range("K74").Value = 3043.4285
$ws.Range("M74").Value = -2169.4285

$ws.Range("H77").Value = 2913
$ws.Range("I77").Value = 3043.4285
$ws.Range("K77").Value = 15217.1425
$ws.Range("M77").Value = -10849.1425

$ws.Range("H136").Value = 1925.1765
$ws.Range("I136").Value = 2565.375
$ws.Range("J136").Value = 1356.1111
$ws.Range("K136").Value = 7696.125
$ws.Range("L136").Value = 4068.3333
$ws.Range("M136").Value = -5146.125
$ws.Range("N136").Value = -9168.3333

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1799.8889
$ws.Range("I134").Value = 1180.8148
$ws.Range("J134").Value = 3657.111
$ws.Range("K134").Value = 3542.4444
$ws.Range("L134").Value = 10971.333
$ws.Range("M134").Value = -1007.4444
$ws.Range("N134").Value = -16041.333

$ws.Range("H140").Value = 41462.5
$ws.Range("J140").Value = 41462.5
$ws.Range("L140").Value = 41462.5
$ws.Range("N140").Value = -51822.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 738.5238000000001
$ws.Range("I7").Value = 1301.7778
$ws.Range("K7").Value = 1301.7778
$ws.Range("M7").Value = -1188.7778

$ws.Range("H22").Value = 2260.2
$ws.Range("I22").Value = 3500.3333
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 3500.3333
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -3150.3333
$ws.Range("N22").Value = -1100

$ws.Range("H132").Value = 5416.9165
$ws.Range("I132").Value = 2670.6667
$ws.Range("J132").Value = 6332.3335
$ws.Range("K132").Value = 8012.000100000001
$ws.Range("L132").Value = 18997.0005
$ws.Range("M132").Value = -5482.000100000001
$ws.Range("N132").Value = -24057.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 520658.3
$ws.Range("J107").Value = 2703112.5
$ws.Range("L107").Value = 8109337.5
$ws.Range("N107").Value = -8113177.5

$ws.Range("H132").Value = 1726.4667
$ws.Range("I132").Value = 2174.75
$ws.Range("J132").Value = 1214.1428
$ws.Range("K132").Value = 19572.75
$ws.Range("L132").Value = 10927.2852
$ws.Range("M132").Value = -17042.75
$ws.Range("N132").Value = -15987.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2132.625
$ws.Range("I102").Value = 1510.1666
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1510.1666
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 111.8334
$ws.Range("N102").Value = -7244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1023
$ws.Range("J22").Value = 1244.4445
$ws.Range("L22").Value = 1244.4445
$ws.Range("N22").Value = -1834.4445

$ws.Range("H27").Value = 1023
$ws.Range("J27").Value = 1244.4445
$ws.Range("L27").Value = 1244.4445
$ws.Range("N27").Value = -1458.4445

$ws.Range("H68").Value = 2731.7856
$ws.Range("I68").Value = 2207
$ws.Range("J68").Value = 3125.375
$ws.Range("K68").Value = 2207
$ws.Range("L68").Value = 3125.375
$ws.Range("M68").Value = -1458
$ws.Range("N68").Value = -4623.375

$ws.Range("H71").Value = 2731.7856
$ws.Range("I71").Value = 2207
$ws.Range("J71").Value = 3125.375
$ws.Range("K71").Value = 11035
$ws.Range("L71").Value = 15626.875
$ws.Range("M71").Value = -7291
$ws.Range("N71").Value = -23114.875

$ws.Range("H123").Value = 40600
$ws.Range("J123").Value = 40600
$ws.Range("L123").Value = 40600
$ws.Range("N123").Value = -50400

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0
